$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '56.959.49'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.83%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.250.10'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '396.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '107.85'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.69%  '
$ws.Range('E7').Value = '  +4.31%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.618'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.24'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.16%  '
$ws.Range('E11').Value = '  +5.62%  '
$ws.Range('E12').Value = '  +2.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.759.44'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.21%  '
$ws.Range('E14').Value = '  +2.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '18.90'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.241.31'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('E17').Value = '  -3.30%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.89'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '56.781.89'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.45%  '
$ws.Range('E20').Value = '  -1.92%  '
$ws.Range('E21').Value = '  +5.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.88'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '291.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.05'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.22%  '
$ws.Range('E25').Value = '  -2.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.04'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '28.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.74%  '
$ws.Range('E28').Value = '  +0.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.30'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.168'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.21%  '
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.110'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.17'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.67%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '40.70'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +11.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0483'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.83%  '
$ws.Range('E36').Value = '  +0.94%  '
$ws.Range('E37').Value = '  -0.46%  '
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('E39').Value = '  -2.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.97'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '137.34'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.22%  '
$ws.Range('E42').Value = '  +1.92%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.283'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.91'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.88%  '
$ws.Range('E45').Value = '  -3.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.60'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '22.11'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.69%  '
$ws.Range('E48').Value = '  +5.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.145.48'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.56%  '
$ws.Range('E50').Value = '  -5.49%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.97'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.04%  '
